$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: English translations ----------------------------------------
# Rows with plain (unstyled / unwrapped) short strings
$ws.Range("C7").Value  = "The railing is broken here..."
$ws.Range("C8").Value  = "Use the rope"
$ws.Range("C9").Value  = "Do nothing"
$ws.Range("C10").Value = "Alright... You can escape from here..."
$ws.Range("C11").Value = "Jump off"

# Rows with long, wrapped dialogue text
$ws.Range("C12").Value = "No, it's way too high up.`nYou would definitely injure yourself.`nIf only you had a rope or something..."
$ws.Range("C13").Value = "No good...`nYou already used the rope to make a hookshot...`nThe hookshot won't be able to reach..."

$ws.Range("C15").Value = "`n<Santa>Good afternoon.`nI'm Santa, the white cat.`nIt's unusual to see a human around here."
$ws.Range("C16").Value = "`n<Santa>Why are you so surprised?`nCats can talk too."
$ws.Range("C17").Value = "`n<Santa>By the way, do you know Ako?`nShe was late coming home, so I got worried and`ncame to check on her. I found all the exits sealed off."
$ws.Range("C18").Value = "`n<Santa>Escape game, hmm.`nSo you're both locked in here.`nThat's a bit of a problem."
$ws.Range("C19").Value = "`n<Santa>...That reminds me.`nI picked up a key down there a little while ago.`nI don't know what it's for."
$ws.Range("C20").Value = "`n<Santa>Well, best of luck escaping.`nIf you manage to get out safe, don't wait up for me.`nPlease give Ako my regards."
$ws.Range("C21").Value = "`n<Santa>Ah, since you're trying to escape, I wouldn't advise`njumping from here. I'm a cat, so I have no problem, but...`nWell... How is your calcium intake?"

# Rows with short strings, but still wrapped / styled
$ws.Range("C23").Value = "Some kind of grate... It's too heavy to open."
$ws.Range("C24").Value = "Lily's Diary"
$ws.Range("C26").Value = "Read"
$ws.Range("C25").Value = "Lily's Diary`nStripes Make a Masochist"

# --- Wrap text on the cells that carry multi-line content -------------------
$ws.Range("B12:B13").WrapText = $true
$ws.Range("B15:B21").WrapText = $true

$ws.Range("C12:C13").WrapText = $true
$ws.Range("C15:C21").WrapText = $true
$ws.Range("C23:C26").WrapText = $true

# --- Row heights for the wrapped rows ---------------------------------------
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 45
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(25).RowHeight = 30

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 48.7109375
$ws.Columns.Item(2).ColumnWidth = 70.5703125
$ws.Columns.Item(3).ColumnWidth = 73.5703125

# --- View state: scroll down and select C25 ----------------------------------
$ws.Range("C25").Select()
$excel.ActiveWindow.ScrollRow = 22

Write-Host "Translations applied"
